# Weekly update: a new price record was inserted at row 300 (Ajo / Feria
# Lagunitas de Puerto Montt), pushing every following record down by one
# row (old row 300 -> new row 301, ..., old row 363 -> new row 364).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 300 (shifts rows 300:363 down to 301:364,
# copying formatting from the row above, same as Excel's native "Insert" UI
# action).
$ws.Rows.Item(300).Insert()

# Populate the newly inserted row 300 with the new weekly record.
$ws.Cells.Item(300, 1).Value  = 4
$ws.Cells.Item(300, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(300, 3).Value  = "Los Lagos"
$ws.Cells.Item(300, 4).Value  = 44889
$ws.Cells.Item(300, 5).Value  = 10
$ws.Cells.Item(300, 6).Value  = 100112003
$ws.Cells.Item(300, 7).Value  = "Ajo"
$ws.Cells.Item(300, 8).Value  = "Chino"
$ws.Cells.Item(300, 9).Value  = "Primera"
$ws.Cells.Item(300, 10).Value = 200
$ws.Cells.Item(300, 11).Value = 18000
$ws.Cells.Item(300, 12).Value = 18000
$ws.Cells.Item(300, 13).Value = 18000
$ws.Cells.Item(300, 14).Value = "`$/caja 10 kilos"
$ws.Cells.Item(300, 15).Value = "China"
$ws.Cells.Item(300, 16).Value = 1800
$ws.Cells.Item(300, 17).Value = 10
$ws.Cells.Item(300, 18).Value = "Hortaliza"
